# Update countries & provincias Spain
# Applies updated case counts and re-sorted country rows to the 'Pais' sheet,
# plus the refreshed "Datos actualizados" timestamp in A1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowUpdates = @(
    ,@(1, @("Datos actualizados a 13 de Abril de 2020 a las 18:22", $null, $null, $null, $null, $null, $null, $null))
    ,@(4, @("Estados Unidos", 564332, 4032, 33735, 507746, 11807, 746, 22851))
    ,@(6, @("Italia", 159516, 3153, 35435, 103616, 3260, 566, 20465))
    ,@(15, @("Suiza", 25688, 273, 12700, 11850, 386, 32, 1138))
    ,@(16, @("Canada", 24833, 450, 7412, 16686, 557, 18, 735))
    ,@(20, @("Austria", 14029, 84, 7343, 6318, 239, 18, 368))
    ,@(32, @("Noruega", 6551, 26, 32, 6386, 59, 5, 133))
    ,@(45, @("Luxemburgo", 3292, 11, 500, 2723, 30, 3, 69))
    ,@(51, @("Singapur", 2918, 386, 586, 2323, 29, 1, 9))
    ,@(52, @("Colombia", 2776, 0, 270, 2397, 92, 0, 109))
    ,@(53, @("Tailandia", 2579, 28, 1288, 1251, 61, 2, 40))
    ,@(58, @("Argelia", 1983, 69, 601, 1069, 60, 20, 313))
    ,@(61, @("Islandia", 1711, 10, 933, 770, 10, 0, 8))
    ,@(79, @("Eslovaquia", 816, 74, 107, 707, 5, 0, 2))
    ,@(80, @("Banglades", 803, 182, 42, 722, 1, 5, 39))
    ,@(109, @("Republica de Yibuti", 298, 84, 41, 255, 0, 0, 2))
    ,@(110, @("Senegal", 291, 11, 178, 111, 1, 0, 2))
    ,@(111, @("Montenegro", 273, 1, 5, 265, 7, 0, 3))
    ,@(112, @("Georgia", 266, 9, 67, 196, 6, 0, 3))
    ,@(113, @("Vietnam", 265, 3, 146, 119, 8, 0, 0))
    ,@(114, @("Guinea", 250, 0, 17, 233, 0, 0, 0))
    ,@(115, @("Isla de Man", 242, 14, 131, 109, 13, 0, 2))
    ,@(116, @("Consejo Danes para los Refugiados", 235, 1, 17, 198, 0, 0, 20))
    ,@(117, @("Sri Lanka", 217, 7, 56, 154, 1, 0, 7))
    ,@(138, @("Togo", 77, 1, 29, 45, 0, 0, 3))
    ,@(163, @("Sudan", 29, 10, 4, 23, 0, 0, 2))
    ,@(164, @("Somalia", 25, 0, 2, 21, 2, 1, 2))
    ,@(165, @("Siria", 25, 0, 5, 18, 0, 0, 2))
    ,@(166, @("Libia", 25, 0, 9, 15, 0, 0, 1))
    ,@(167, @("Antigua y Barbuda", 23, 2, 0, 21, 1, 0, 2))
    ,@(168, @("Republica del Chad", 23, 5, 2, 21, 0, 0, 0))
    ,@(169, @("Mozambique", 21, 0, 2, 19, 0, 0, 0))
    ,@(170, @("Guinea Ecuatorial", 21, 0, 3, 18, 0, 0, 0))
    ,@(171, @("Maldivas", 20, 0, 14, 6, 0, 0, 0))
    ,@(172, @("Laos", 19, 0, 0, 19, 0, 0, 0))
    ,@(173, @("Angola", 19, 0, 4, 13, 0, 0, 2))
    ,@(174, @("Nueva Caledonia", 18, 0, 1, 17, 0, 0, 0))
    ,@(205, @("Burundi", 5, 0, 0, 4, 0, 1, 1))
    ,@(206, @("Islas Malvinas", 5, 0, 1, 4, 0, 0, 0))
    ,@(215, @("Yemen", 1, 0, 0, 1, 0, 0, 0))
    ,@(216, @("San Pedro y Miquelon", 1, 0, 0, 1, 0, 0, 0))
)

foreach ($entry in $rowUpdates) {
    $r = $entry[0]
    $vals = $entry[1]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        if ($null -ne $vals[$i]) {
            $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
        }
    }
}

